$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2 = 0.16875
    3 = 0.1213315812527376
    4 = 0.1122727272727273
    5 = 0.1540084388185654
    6 = 0.157258064516129
    7 = 0.211890243902439
    8 = 0.09012875536480687
    9 = 0.1275964391691395
    10 = 0.1018387553041018
    11 = 0.1124620060790274
    12 = 0.1566731141199226
    13 = 0.08431372549019608
    14 = 0.07644787644787644
    15 = 0.2054054054054054
    16 = 0.1467089611419508
    17 = 0.1666666666666667
    18 = 0.1315192743764172
    19 = 0.1003344481605351
    20 = 0.125
    21 = 0.08487084870848709
    22 = 0.07977207977207977
    23 = 0.1220132180986273
    24 = 0.1495327102803738
    25 = 0.1035353535353535
    26 = 0.09826589595375723
    27 = 0.1117021276595745
    28 = 0.07230769230769231
    29 = 0.1400709219858156
    30 = 0.1122448979591837
    31 = 0.154320987654321
    32 = 0.1155778894472362
    33 = 0.1596385542168675
    34 = 0.1804812834224599
    35 = 0.1838842975206612
    36 = 0.1360153256704981
    37 = 0.03683241252302026
    38 = 0.1047227926078029
    39 = 0.08450704225352113
    40 = 0.1256544502617801
    41 = 0.1328358208955224
    42 = 0.1085271317829457
    43 = 0.05636363636363637
    44 = 0.1189336978810663
    45 = 0.1122194513715711
    46 = 0.1242215666994428
    47 = 0.09919028340080972
    48 = 0.08983957219251337
    49 = 0.191588785046729
    50 = 0.1071800208116545
    51 = 0.08765315739868049
    52 = 0.1233644859813084
    53 = 0.0970464135021097
    54 = 0.1323251417769376
    55 = 0.118964310706788
    56 = 0.103448275862069
    57 = 0.08975521305530372
    58 = 0.1147859922178988
    59 = 0.1631578947368421
    60 = 0.1098901098901099
    61 = 0.1400523560209424
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 12).Value = $values[$row]
}
